$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Scratch cell, far away from any real data, used purely to force numeric-looking
# strings ("349467", "075026", ...) to be stored as text (shared string) instead
# of being auto-coerced to numbers. It is cleared again before the script ends.
$scratch = $ws1.Range("ZZ1")

function Set-TextValue($range, $text) {
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

# --- Apply the new randomized "Customer ID" values (kept as text) ---
Set-TextValue $ws1.Range("F2") "349467"
Set-TextValue $ws1.Range("F3") "356300"
Set-TextValue $ws1.Range("F4") "584014"
Set-TextValue $ws1.Range("F5") "695146"
Set-TextValue $ws1.Range("F6") "969209"
Set-TextValue $ws1.Range("F7") "253850"
Set-TextValue $ws1.Range("F8") "075026"

$ws1.Application.CutCopyMode = $false
$scratch.Clear()

# --- Move "Customer ID" column data (rows 4-8) from sheet1 to a new column F on sheet2 ---
$ws1.Range("F2").Copy()
$ws2.Range("F1:F5").PasteSpecial(-4122)

$ws1.Range("F4:F8").Copy()
$ws2.Range("F1:F5").PasteSpecial(-4163)

$ws1.Application.CutCopyMode = $false

# --- Remove rows 4-8 (and stale hyperlinks) from sheet1 now that the data lives on sheet2 ---
$ws1.Range("A4:F8").EntireRow.Delete()
$ws1.Hyperlinks.Delete()

# --- Rename the "Customer ID" header to "Customer_ID" ---
$ws1.Range("F1").Value = "Customer_ID"

# --- Update the active selection shown on each sheet ---
$ws2.Range("C9").Select()
$ws1.Range("D20").Select()
$ws1.Activate()
